# Reorder the data rows (2-7) on Sheet1: rows were re-sorted by
# Profitcenter/Belegnummer priority ("refactoring bpmn components fixed
# flow prioritiese"). Apply the new row values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-7, columns A (Belegnummer), B (Belegdatum),
# C (Lieferant), D (Profitcenter), E (Kostenstelle)
$data = @(
    @(1, "05.Jan", "A", "P1", "K2"),
    @(3, "07.Jan", "A", "P1", "K2"),
    @(5, "09.Jan", "B", "P2", "K1"),
    @(7, "11.Jan", "A", "P2", "K1"),
    @(2, "06.Jan", "B", "P3", "unknown"),
    @(6, "10.Jan", "A", "P3", "unknown")
)

$rowIndex = 2
foreach ($rowValues in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rowValues[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rowValues[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rowValues[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rowValues[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rowValues[4]
    $rowIndex++
}
